# Add a new worksheet "nr_studies" at the end of the workbook, reporting
# the number of effect sizes and number of studies (k) for each
# outcome x moderator_context combination.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab strip (matches sheetId="4" / rId4 ordering in the target file).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "nr_studies"

# Header row
$ws.Range("A1").Value = "outcome"
$ws.Range("B1").Value = "moderator_context"
$ws.Range("C1").Value = "n_effect_sizes"
$ws.Range("D1").Value = "k_studies"
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108

# Data rows: outcome, moderator_context, n_effect_sizes, k_studies
$data = @(
    @("NS","Education",355,33),
    @("NS","Parenting",186,27),
    @("NS","Sport",127,13),
    @("NS","Exercise",37,5),
    @("NS","Healthcare",4,1),
    @("NS","Health care",10,1),
    @("NS","Partner",3,1),
    @("NT","Education",109,15),
    @("NT","Parenting",107,17),
    @("NT","Sport",111,12),
    @("NT","Exercise",44,6),
    @("NT","Health care",10,1)
)

$rowIdx = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIdx, 1).Value = $row[0]
    $ws.Cells.Item($rowIdx, 2).Value = $row[1]
    $ws.Cells.Item($rowIdx, 3).Value = $row[2]
    $ws.Cells.Item($rowIdx, 4).Value = $row[3]
    $rowIdx++
}

Write-Output "Added worksheet 'nr_studies' with $($data.Count) data rows."
